# Rotate the taxon-record data across rows 34,35,36,37,39,40 as described
# by the upstream diff. Column C (Valideringsstatus) and the other
# per-row metadata (P, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AW, AX) stay
# untouched - only A,B,D,E,F,G,H,Q,R change values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    34 = @{ A = 112038601; B = 73634;  D = "LC"; E = 6426; F = "Kattfotslav";     G = "Felipes leucopellaeus";   H = "(Ach.) Frisch & G.Thor";  Q = 616012.5978259755; R = 6895611.944218947 }
    35 = @{ A = 112038602; B = 86223;  D = "NT"; E = 4412; F = "Äggvaxskivling";  G = "Hygrophorus karstenii";   H = "Sacc. & Cub.";            Q = 616026.2967975155; R = 6895553.979090866 }
    36 = @{ A = 112038596; B = 90087;  D = "LC"; E = 3298; F = "Trådticka";       G = "Climacocystis borealis";  H = "(Fr.) Kotl. & Pouzar";    Q = 616076.0611235843; R = 6895427.595461337 }
    37 = @{ A = 112038600; B = 86223;  D = "NT"; E = 4412; F = "Äggvaxskivling";  G = "Hygrophorus karstenii";   H = "Sacc. & Cub.";            Q = 616034.1211971109; R = 6895585.10294092  }
    39 = @{ A = 112038603; B = 89369;  D = "LC"; E = 5447; F = "Vedticka";        G = "Fuscoporia viticola";     H = "(Schwein.) Murrill";      Q = 615968.1934313668; R = 6895405.650930508 }
    40 = @{ A = 112038604; B = 89845;  D = "VU"; E = 1209; F = "Rynkskinn";       G = "Phlebia centrifuga";      H = "P.Karst.";                Q = 615977.7276359925; R = 6895550.438170813 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("A$r").Value = $vals.A
    $ws.Range("B$r").Value = $vals.B
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("E$r").Value = $vals.E
    $ws.Range("F$r").Value = $vals.F
    $ws.Range("G$r").Value = $vals.G
    $ws.Range("H$r").Value = $vals.H
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("R$r").Value = $vals.R
}
